# Add a new "type" (grammar category) column (K) to the hidden "grammar_2"
# worksheet, populating it with the Korean grammar-topic label that
# corresponds to each existing row, and update the sheet's selection.

$wb = $excel.ActiveWorkbook

# The data lives on the hidden "grammar_2" sheet (sheet1.xml); "Sheet1" is
# the workbook's visible/active tab and must remain so.
$ws       = $wb.Worksheets.Item("grammar_2")
$visibleWs = $wb.Worksheets.Item("Sheet1")

# New header + per-row category labels for column K (rows 1-51).
$values = @(
    "type",
    "수량형용사", "수량형용사", "수량형용사", "수량형용사", "수량형용사", "수량형용사",
    "빈도부사", "빈도부사", "빈도부사", "빈도부사", "빈도부사", "빈도부사",
    "현재진행형", "현재진행형", "현재진행형", "현재진행형", "현재진행형", "현재진행형",
    "부정대명사와 수량형용사", "부정대명사와 수량형용사", "부정대명사와 수량형용사", "부정대명사와 수량형용사", "부정대명사와 수량형용사", "부정대명사와 수량형용사", "부정대명사와 수량형용사",
    "조동사", "조동사", "조동사", "조동사", "조동사",
    "미래시제", "미래시제", "미래시제", "미래시제", "미래시제", "미래시제",
    "과거시제", "과거시제", "과거시제", "과거시제", "과거시제", "과거시제", "과거시제", "과거시제",
    "의문문", "의문문", "의문문", "의문문", "의문문", "의문문"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 11).Value = $values[$i]
}

# Reflect the new selection/view on the grammar_2 sheet (activate it so the
# selection is recorded, then restore Sheet1 as the active tab so the
# workbook-level active-sheet state is unchanged).
$ws.Activate()
$ws.Range("K46").Select()
$visibleWs.Activate()
